$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto-price refresh diff.
# Values that look like plain decimal numbers (e.g. "1.00", "26.00")
# are forced to stay text (matching the source inlineStr cells) by
# temporarily marking the cell as Text-formatted, then restoring the
# default "Normal" style so no new formatting is left behind.

$ws.Range("D2").Value = "46.000.55"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "2.449.02"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +5.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "2.822.34"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "2.432.84"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "45.858.62"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "0.0₃0932"
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -4.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("E32").Value = "  +6.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0761"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.24%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "1.963.46"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.74%  "
